$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-affecting data: rows 2-9 get new values, plus a brand new row 10 (Q8).

# Row 2 (Q0)
$ws.Range("B2").Value = -0.0409949888987393
$ws.Range("C2").Value = 0.1649705476528477
$ws.Range("D2").Value = 0.05860056171370987
$ws.Range("E2").Value = 0.2420755289443976
$ws.Range("F2").Value = 0.2475852015561247

# Row 3 (Q1)
$ws.Range("B3").Value = -0.0582393705973169
$ws.Range("C3").Value = 0.2339998014479659
$ws.Range("D3").Value = 0.1090099205812892
$ws.Range("E3").Value = 0.3301665043296931
$ws.Range("F3").Value = 0.33825966995519

# Row 4 (Q2)
$ws.Range("B4").Value = -0.05540528605858396
$ws.Range("C4").Value = 0.2730276854007799
$ws.Range("D4").Value = 0.1281020150765566
$ws.Range("E4").Value = 0.357913418408079
$ws.Range("F4").Value = 0.3693221348559199

# Row 5 (Q3)
$ws.Range("B5").Value = 0.0007705445796656607
$ws.Range("C5").Value = 0.1970834307540707
$ws.Range("D5").Value = 0.06344748247054742
$ws.Range("E5").Value = 0.251887837083388
$ws.Range("F5").Value = 0.2641809561735251

# Row 6 (Q4)
$ws.Range("B6").Value = 0.03211500412566347
$ws.Range("C6").Value = 0.1962129447978213
$ws.Range("D6").Value = 0.05653652986779405
$ws.Range("E6").Value = 0.237774115218192
$ws.Range("F6").Value = 0.2483392759418782

# Row 7 (Q5)
$ws.Range("B7").Value = 0.01002085364717452
$ws.Range("C7").Value = 0.1353809626663662
$ws.Range("D7").Value = 0.02435385723373847
$ws.Range("E7").Value = 0.1560572242279686
$ws.Range("F7").Value = 0.1651820804193373
$ws.Range("G7").Value = 9

# Row 8 (Q6)
$ws.Range("B8").Value = -0.02479498802512995
$ws.Range("C8").Value = 0.2038757822666403
$ws.Range("D8").Value = 0.05563351498145874
$ws.Range("E8").Value = 0.2358675793352252
$ws.Range("F8").Value = 0.2569483766447083
$ws.Range("G8").Value = 6

# Row 9 (Q7)
$ws.Range("B9").Value = 0.1034004173967734
$ws.Range("C9").Value = 0.1076593290439133
$ws.Range("D9").Value = 0.0258513414274284
$ws.Range("E9").Value = 0.1607835234948793
$ws.Range("F9").Value = 0.1507963615754776
$ws.Range("G9").Value = 3

# New row 10 (Q8)
$ws.Range("A10").Value = "Q8"
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").Value = -0.05594074298689113
$ws.Range("C10").Value = 0.05594074298689113
$ws.Range("D10").Value = 0.003129366725925409
$ws.Range("E10").Value = 0.05594074298689113
$ws.Range("G10").Value = 1
